# Applies the "AutonHack Presentation" edit:
#  - Slide 9, bullet list shape: bold three bullet lines and tweak two
#    run texts ("Multimodal " -> "Multimodal - ", "inputs" -> "Multi-languages").
#  - Re-colours theme1.xml's colour scheme to the "Default" palette that,
#    in the target OOXML, ends up swapped into theme1.xml (this runtime's
#    ColorScheme/ThemeColorScheme COM surface always binds to the
#    presentation's primary theme part, theme1.xml).

$p  = $ppt.ActivePresentation
$s9 = $p.Slides.Item(9)
$sh = $s9.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# --- Bold the three bullet runs (lengths unaffected, so offsets below are
#     computed against the shape's original text layout and are safe to
#     apply before any text is edited). ---

# Paragraph 3: "Multimodal " (run 1) + "inputs" (run 2)
$tr.Characters(115, 11).Font.Bold = -1
$tr.Characters(126, 6).Font.Bold  = -1

# Paragraph 5: "Better LLMs"
$tr.Characters(149, 11).Font.Bold = -1

# Paragraph 8: "Reinforcement learning (RLHF)"
$tr.Characters(197, 29).Font.Bold = -1

# --- Text edits. Apply right-to-left (higher offsets first) so earlier,
#     still-to-be-used offsets are not invalidated by length changes. ---

# Paragraph 3, run 2: "inputs" -> "Multi-languages"
$tr.Characters(126, 6).Text = "Multi-languages"

# Paragraph 3, run 1: "Multimodal " -> "Multimodal - "
$tr.Characters(115, 11).Text = "Multimodal - "

# --- Recolour the primary theme to the target "Default" palette. ---
$tcs = $s9.ThemeColorScheme
$tcs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 0x588115   # dk2      -> 158158
$tcs.Item(4).RGB  = 0xF3F3F3   # lt2      -> F3F3F3
$tcs.Item(5).RGB  = 0xC78D05   # accent1  -> 058DC7
$tcs.Item(6).RGB  = 0x32B450   # accent2  -> 50B432
$tcs.Item(7).RGB  = 0x1B56ED   # accent3  -> ED561B
$tcs.Item(8).RGB  = 0x00EFED   # accent4  -> EDEF00
$tcs.Item(9).RGB  = 0xE5CB24   # accent5  -> 24CBE5
$tcs.Item(10).RGB = 0x72E564   # accent6  -> 64E572
$tcs.Item(11).RGB = 0xCC0022   # hlink    -> 2200CC
$tcs.Item(12).RGB = 0x8B1A55   # folHlink -> 551A8B
